# Simulated Wild Card round and logged it.
# Appends this game's play-by-play yardage samples to the season-long
# shared-string lists on YDS/ST, and rolls the per-game totals forward
# on OFF / DEF / ST / TURNS / PEN.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# YDS sheet: append this game's individual play yardages to the
# season-long space-separated lists.
# ---------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 3 6 7 -1 8 4 15 -1 3 5 1 7 12 9 2 4 6 19 0 4 4 6 5 -1 2 4 1 2 1 16 5 2 2 3 2 3 14 14 0 0 1 8 2 7 1 5 1"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 4 2 5 4 3 0 7 5 2 1 10 6 5 3 1 3 3 11 9 1 27 7 3 5 1 3 1 4 3 2 2 23 4 -1 5 2 -2 4 11 3 2 0 4 26 15 1 7 6 5 9 3 8 2 4 9 1 3 9 16 2 15 10 1 3 7 6 4 0"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 2 5 35 7 5 13 23 17 6 2 28 6 32 5 1 18 39 11 8 9 30 4 12 10 2 5 8 7 10 7 1 43 4 7 15 3 6 2 18"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 9 4 5 13 7 9 7 17 4 3 5 9 2 0 15 9 0 7 8 24 2 19 6 11 22 19 45 4 9 4 34 8 19 19 38 1"

# ---------------------------------------------------------------
# OFF sheet: roll this game's offensive totals into the season sums.
# ---------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("B2").Value = 10
$offWs.Range("C2").Value = 448
$offWs.Range("D2").Value = 25
$offWs.Range("F2").Value = 130
$offWs.Range("G2").Value = 141
$offWs.Range("I2").Value = 14
$offWs.Range("J2").Value = 88
$offWs.Range("N2").Value = 33
$offWs.Range("O2").Value = 37
$offWs.Range("P2").Value = 24

$offWs.Range("C3").Value = 315
$offWs.Range("D3").Value = 8
$offWs.Range("E3").Value = 69
$offWs.Range("F3").Value = 194
$offWs.Range("G3").Value = 54
$offWs.Range("H3").Value = 42
$offWs.Range("I3").Value = 104
$offWs.Range("J3").Value = 96
$offWs.Range("L3").Value = 505
$offWs.Range("M3").Value = 332
$offWs.Range("Q3").Value = 1062

# ---------------------------------------------------------------
# DEF sheet: roll this game's defensive totals into the season sums.
# ---------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value = 420
$defWs.Range("D2").Value = 27
$defWs.Range("E2").Value = 20
$defWs.Range("F2").Value = 133
$defWs.Range("G2").Value = 143
$defWs.Range("I2").Value = 18
$defWs.Range("J2").Value = 81
$defWs.Range("N2").Value = 30
$defWs.Range("O2").Value = 48

$defWs.Range("B3").Value = 20
$defWs.Range("C3").Value = 353
$defWs.Range("F3").Value = 203
$defWs.Range("G3").Value = 78
$defWs.Range("H3").Value = 49
$defWs.Range("I3").Value = 99
$defWs.Range("J3").Value = 102
$defWs.Range("L3").Value = 490
$defWs.Range("M3").Value = 311
$defWs.Range("Q3").Value = 984

# ---------------------------------------------------------------
# ST sheet: roll this game's special-teams totals into the season
# sums, and append the game's individual kick/punt distances to the
# season-long lists.
# ---------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 179
$stWs.Range("D2").Value = 107
$stWs.Range("F2").Value = 440
$stWs.Range("G2").Value = 428
$stWs.Range("L2").Value = 122
$stWs.Range("M2").Value = 91
$stWs.Range("B3").Value = 85

$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 46 45 47 47 42 53"
$stWs.Range("B4").Value = $stWs.Range("B4").Value2 + " 62 62 50 58 57 42"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 0 0 0 0 2 52"
$stWs.Range("B5").Value = $stWs.Range("B5").Value2 + " 20 18 15 5 15 0"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 11 0 9 0"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 24 15 16 21 25 26"

# ---------------------------------------------------------------
# TURNS sheet: roll this game's turnover totals into the season sums.
# ---------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("B3").Value = 13

# ---------------------------------------------------------------
# PEN sheet: roll this game's penalty totals into the season sums.
# ---------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("B2").Value = 23
$penWs.Range("D2").Value = 18
$penWs.Range("D4").Value = 15
